$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update existing discount amounts
$ws.Range("D18").Value = 3
$ws.Range("D30").Value = 15
$ws.Range("D31").Value = 7

# Append a new row (row 47) with a new discount entry
$ws.Range("A47").Value = 46
$ws.Range("B47").Value = 7
$ws.Range("C47").Value = 48
$ws.Range("D47").Value = 4
$ws.Range("E47").Value = "System"
$ws.Range("F47").Value = "2025-03-16 19:02:52"
$ws.Range("G47").Value = 0
